# Fixed the I/O explanation
# Slide 5, "Content Placeholder 2" (Shape 2) contains the bullet list with
# the File/Socket I/O bullets that need correcting.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Edit 1 ---------------------------------------------------------------
# "Provides comprehensive File and Socket I/O through Foundation library."
#   -> "...through predefined set of classes/APIs."
# This paragraph only has a single run, so just rewrite its text in place.
$para14 = $tr.Paragraphs(14, 1)
$run14  = $para14.Runs(1, 1)
$run14.Text = "Provides comprehensive File and Socket I/O through predefined set of classes/APIs."

# --- Edit 2 ---------------------------------------------------------------
# "File I/O classes: ... NSInputStream etc." gains a trailing clause:
# "... NSInputStream etc from Foundation library." split up as
# " " / "etc" / " " / "from " / "Foundation library." runs.
$para15 = $tr.Paragraphs(15, 1)

# The last run of the paragraph is " etc."; swap its trailing "." character
# for the new " from Foundation library." tail (keeps "etc" itself intact,
# drops the period that used to end the sentence).
$tail = $para15.Characters(89, 1)
$tail.Text = " from Foundation library."

# Re-touch each new word/space with its (unchanged) font size so the host
# engine materialises them as their own runs instead of one big run.
$etcRun = $para15.Characters(86, 3)
$etcRun.Font.Size = 18

$fromRun = $para15.Characters(90, 5)
$fromRun.Font.Size = 18

$foundationRun = $para15.Characters(95, 19)
$foundationRun.Font.Size = 18

Write-Host "Para14:" $para14.Text
Write-Host "Para15:" $para15.Text
